$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last status check" timestamp in F1
$ws.Range("F1").Value = "Last status check on: 23.02.2022 20:45"

# Shift the price history for row 10 (EuroOil Opustena):
# old "new price" (B10) becomes "old price" (C10)
$ws.Range("C10").Value = $ws.Range("B10").Value2

# set the newly observed price
$ws.Range("B10").Value = 37.7

# set the delta as a formatted text string (leading apostrophe forces text,
# preventing Excel from re-parsing "+0.3" back into a number)
$ws.Range("D10").Value = "'+0.3"
$ws.Range("D10").Style = "Normal"

# set the observed timestamp as plain text, clearing the previous date style
$ws.Range("E10").Value = "2022-02-23 20:45:32"
$ws.Range("E10").Style = "Normal"
